# Append two more days of GSC export data to the "Chart" sheet (rows 70-71),
# mirroring the pattern of the existing rows:
#   column A = date label (stored as text, like the existing rows)
#   column B = Invalid count
#   column C = Valid count

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 70: 2025-12-13
# Write the date as a formula that evaluates to a text string, then convert it
# to a literal value via copy/paste-special. This avoids Excel's automatic
# "looks like a date" conversion (which would store it as a date serial
# number instead of text, like the rest of column A).
$ws.Range("A70").Formula = '="2025-12-13"'
$ws.Range("A70").Copy()
$ws.Range("A70").PasteSpecial(-4163)
$ws.Range("B70").Value = 0
$ws.Range("C70").Value = 30

# Row 71: 2025-12-14
$ws.Range("A71").Formula = '="2025-12-14"'
$ws.Range("A71").Copy()
$ws.Range("A71").PasteSpecial(-4163)
$ws.Range("B71").Value = 0
$ws.Range("C71").Value = 31

Write-Output "Added rows 70-71 to Chart sheet"
